$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for row 5 (October 4th work entry), following the same column
# layout already used by the existing rows 2-4.
$ws.Range("B5").Value = "Enterainer"
$ws.Range("C5").Value = "Énergie"
$ws.Range("E5").Value = "Anxiété, Aveugle"
$ws.Range("F5").Value = "Champs de bataille"
$ws.Range("G5").Value = "Arts martiaux"
$ws.Range("H5").Value = "Culbutes"
$ws.Range("I5").Value = "Armes de jet"
$ws.Range("J5").Value = "Ménestrel"

# Cell styles mirror the conditional "Good/Neutral/Bad" pattern already
# used in rows 2-4; copy the format from a cell in the same column that
# already carries the right look so the workbook's existing style table
# (with its border + centered alignment) is reused instead of rebuilt.
$xlPasteFormats = -4122

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B5").PasteSpecial($xlPasteFormats)

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C5").PasteSpecial($xlPasteFormats)

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D5").PasteSpecial($xlPasteFormats)

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E5").PasteSpecial($xlPasteFormats)

$ws.Range("F3").Copy() | Out-Null
$ws.Range("F5").PasteSpecial($xlPasteFormats)

$ws.Range("G2").Copy() | Out-Null
$ws.Range("G5").PasteSpecial($xlPasteFormats)

$ws.Range("H2").Copy() | Out-Null
$ws.Range("H5").PasteSpecial($xlPasteFormats)

$ws.Range("I2").Copy() | Out-Null
$ws.Range("I5").PasteSpecial($xlPasteFormats)

$ws.Range("J2").Copy() | Out-Null
$ws.Range("J5").PasteSpecial($xlPasteFormats)

$ws.Range("K2").Copy() | Out-Null
$ws.Range("K5").PasteSpecial($xlPasteFormats)

$ws.Range("L2").Copy() | Out-Null
$ws.Range("L5").PasteSpecial($xlPasteFormats)

$ws.Range("M2").Copy() | Out-Null
$ws.Range("M5").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
